# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values, plus the EnergySwap/RenderToken row swap (rows 45-46, columns B-E).
#
# New cell values are applied as literal text. Several "Price" strings are
# valid-looking numbers (e.g. "1.002", "22.25"); Excel's COM layer will
# silently coerce a bare numeric-looking string to a real number when it is
# assigned to a General-formatted cell, which would corrupt both the stored
# type and precision (e.g. "1.001" -> 1.0009999999999999). To keep those
# cells as plain text - matching the original inline-string cells - the
# target cell is temporarily switched to the Text number format ("@") before
# the value is written, then its style is restored from the sibling "Coin"
# cell in the same row (column B), which has always carried the workbook's
# default (unstyled) format. This keeps formatting identical to the
# original file while forcing the value to be stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for every changed cell (A1-style ref -> new text).
$updates = [ordered]@{
    "D2" = "29.134.25"
    "E2" = "  +1.27%  "
    "D3" = "1.905.18"
    "E3" = "  +1.59%  "
    "D4" = "1.002"
    "E4" = "  -0.17%  "
    "D5" = "326.92"
    "E5" = "  +0.65%  "
    "E6" = "  -0.18%  "
    "D7" = "0.4605"
    "E7" = "  +0.14%  "
    "D8" = "0.3935"
    "E8" = "  +1.60%  "
    "D9" = "46.66"
    "E9" = "  +1.03%  "
    "D10" = "0.07928"
    "D11" = "1.001"
    "E11" = "  +1.16%  "
    "D12" = "22.25"
    "E12" = "  +2.00%  "
    "D13" = "1.894.66"
    "E13" = "  +2.02%  "
    "D14" = "7.087"
    "E14" = "  +1.25%  "
    "D15" = "5.754"
    "E15" = "  +0.71%  "
    "D16" = "0.06942"
    "E16" = "  -0.35%  "
    "D17" = "88.31"
    "E17" = "  -0.19%  "
    "E18" = "  -0.22%  "
    "E19" = "  +0.10%  "
    "D20" = "17.09"
    "E20" = "  +1.81%  "
    "E21" = "  -0.12%  "
    "D22" = "29.154.62"
    "E22" = "  +1.35%  "
    "D23" = "5.359"
    "E23" = "  +1.37%  "
    "E24" = "  +0.35%  "
    "D25" = "2.136.54"
    "E25" = "  +2.53%  "
    "D26" = "2.055"
    "E26" = "  -2.14%  "
    "D27" = "156.52"
    "E27" = "  +2.28%  "
    "D28" = "19.42"
    "E28" = "  +0.99%  "
    "D29" = "6.119"
    "E29" = "  +4.34%  "
    "D30" = "1.995"
    "E30" = "  +0.91%  "
    "D31" = "118.62"
    "E31" = "  -0.40%  "
    "D32" = "0.09376"
    "D33" = "0.9291"
    "D34" = "5.327"
    "E34" = "  +0.27%  "
    "D35" = "1.349"
    "E35" = "  +0.70%  "
    "D36" = "3.268"
    "E36" = "  -1.66%  "
    "D37" = "1.210"
    "E37" = "  +5.17%  "
    "D38" = "0.05829"
    "E38" = "  +0.87%  "
    "D39" = "0.02102"
    "E39" = "  +1.32%  "
    "D40" = "7.920"
    "E40" = "  +2.96%  "
    "E41" = "  -0.22%  "
    "D42" = "0.5741"
    "E42" = "  +1.63%  "
    "D43" = "0.1800"
    "E43" = "  +0.72%  "
    "D44" = "9.946"
    "E44" = "  +0.83%  "
    "B45" = "RenderToken"
    "C45" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D45" = "2.240"
    "E45" = "  +4.31%  "
    "B46" = "EnergySwap"
    "C46" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D46" = "11.92"
    "E46" = "  +1.09%  "
    "D47" = "0.5411"
    "E47" = "  +2.18%  "
    "D48" = "0.07069"
    "E48" = "  -2.02%  "
    "E49" = "  +2.42%  "
    "D50" = "2.556"
    "E50" = "  +5.98%  "
    "D51" = "113.07"
    "E51" = "  -0.38%  "
}

# Cells whose new text would otherwise be auto-converted to a number.
$forceText = @(
    "D4","D5","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17",
    "D20","D23","D26","D27","D28","D29","D30","D31","D32","D33","D34",
    "D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46",
    "D47","D48","D50","D51"
)

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $newValue = $updates[$ref]

    if ($forceText -contains $ref) {
        $row = $ref.Substring(1)
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = $ws.Range("B$row").Style
    } else {
        $cell.Value = $newValue
    }
}
